$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data

$ws.Range("D2").Value = "'" + '29.804.50'
$ws.Range("E2").Value = '  +1.55%  '

$ws.Range("D3").Value = "'" + '1.933.71'
$ws.Range("E3").Value = '  +1.28%  '

$ws.Range("D4").Value = "'" + '1.009'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = "'" + '337.61'
$ws.Range("E5").Value = '  +4.08%  '

$ws.Range("D6").Value = "'" + '1.007'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").Value = "'" + '0.4832'
$ws.Range("E7").Value = '  +0.46%  '

$ws.Range("D8").Value = "'" + '0.4107'
$ws.Range("E8").Value = '  +1.14%  '

$ws.Range("D9").Value = "'" + '0.08160'
$ws.Range("E9").Value = '  -0.51%  '

$ws.Range("D10").Value = "'" + '1.013'
$ws.Range("E10").Value = '  -0.66%  '

$ws.Range("D11").Value = "'" + '23.66'
$ws.Range("E11").Value = '  +0.75%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'" + '1.987.63'
$ws.Range("E12").Value = '  +2.59%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'" + '6.070'
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("D14").Value = "'" + '7.269'
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").Value = "'" + '90.80'

$ws.Range("D16").Value = "'" + '0.06839'
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("D18").Value = "'" + '0.00001033'
$ws.Range("E18").Value = '  -0.25%  '

$ws.Range("D19").Value = "'" + '17.75'
$ws.Range("E19").Value = '  +0.41%  '

$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").Value = "'" + '29.788.54'
$ws.Range("E21").Value = '  +1.36%  '

$ws.Range("D22").Value = "'" + '5.626'
$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("E23").Value = '  +0.62%  '

$ws.Range("D24").Value = "'" + '2.177'
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("D25").Value = "'" + '2.133.26'
$ws.Range("E25").Value = '  -0.72%  '

$ws.Range("D26").Value = "'" + '6.570'
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").Value = "'" + '156.98'
$ws.Range("E27").Value = '  +0.49%  '

$ws.Range("D28").Value = "'" + '20.02'

$ws.Range("D29").Value = "'" + '2.088'
$ws.Range("E29").Value = '  -0.72%  '

$ws.Range("D30").Value = "'" + '120.90'
$ws.Range("E30").Value = '  +0.61%  '

$ws.Range("D31").Value = "'" + '1.006'
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("D32").Value = "'" + '0.09654'
$ws.Range("E32").Value = '  +1.17%  '

$ws.Range("D33").Value = "'" + '5.537'
$ws.Range("E33").Value = '  -0.83%  '

$ws.Range("D34").Value = "'" + '1.408'
$ws.Range("E34").Value = '  +3.31%  '

$ws.Range("D35").Value = "'" + '3.532'
$ws.Range("E35").Value = '  -0.39%  '

$ws.Range("D36").Value = "'" + '0.06580'
$ws.Range("E36").Value = '  +7.77%  '

$ws.Range("D37").Value = "'" + '0.02280'
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").Value = "'" + '1.199'
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D39").Value = "'" + '0.5965'
$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("D40").Value = "'" + '10.76'
$ws.Range("E40").Value = '  -0.39%  '

$ws.Range("D41").Value = "'" + '7.925'
$ws.Range("E41").Value = '  -1.48%  '

$ws.Range("D42").Value = "'" + '0.1846'
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("D43").Value = "'" + '2.470'
$ws.Range("E43").Value = '  +3.47%  '

$ws.Range("D44").Value = "'" + '1.273'
$ws.Range("E44").Value = '  -0.36%  '

$ws.Range("D45").Value = "'" + '12.26'
$ws.Range("E45").Value = '  -1.35%  '

$ws.Range("D46").Value = "'" + '0.07475'
$ws.Range("E46").Value = '  -1.71%  '

$ws.Range("D47").Value = "'" + '0.5550'
$ws.Range("E47").Value = '  -0.36%  '

$ws.Range("D48").Value = "'" + '1.981'
$ws.Range("E48").Value = '  +1.61%  '

$ws.Range("D49").Value = "'" + '116.76'

$ws.Range("D50").Value = "'" + '2.414'
$ws.Range("E50").Value = '  -0.29%  '

$ws.Range("D51").Value = "'" + '72.24'
$ws.Range("E51").Value = '  +0.20%  '
